$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test")

# Update the EActorType enum value labels to their fully-qualified names
$ws.Range("C5").Value = "ACTOR_TYPE_NONE"
$ws.Range("C6").Value = "ACTOR_TYPE_PLAYER"
$ws.Range("C7").Value = "ACTOR_TYPE_NPC"

# Column C widened (best-fit) to accommodate the new, longer enum names
$ws.Columns.Item(3).ColumnWidth = 17.285714285714285

# Move the active selection
$ws.Range("D10").Select()
